$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 5, 6 and 7 (old entries for records 3, 4 and 5)
$ws.Rows("5:7").Delete()

# Update row 4 with new data (record replaced with a different process)
$ws.Range("B4").Value = "5002019-61.2019.8.21.0067"
$ws.Range("C4").Value = "0064471-68.2019.8.21.9000"
$ws.Range("D4").Value = "Relacionado na TR"
$ws.Range("E4").Value = "0042940-86.2020.8.21.9000"
$ws.Range("F4").Value = "Relacionado na TR"
$ws.Range("G4").Value = "9000401-13.2019.8.21.0067"
$ws.Range("H4").Value = "Migrado"
$ws.Range("I4").Value = "27/05/2019"
